# BasicApp: use kobalt email address for admin user
#
# The "TechnicalUser" sheet holds a single bootstrap admin account whose
# firstName/lastName used to be the single French word "Administrateur" and
# whose e-mail was the placeholder admin@example.com. This change splits the
# name into firstName "Admin" / lastName "Kobalt" and switches the e-mail
# address (cell value, and hyperlink address/display text) to the Kobalt
# qualification inbox. It also leaves the TechnicalUser sheet selected/active
# instead of BasicUser.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("TechnicalUser")

# --- Row 2 (the admin account) -------------------------------------------
$ws3.Range("B2").Value = "Admin"
$ws3.Range("C2").Value = "Kobalt"
$ws3.Range("E2").Value = "interne-qualif@kobalt.fr"

# --- Hyperlink on the e-mail cell: point it at the new address and show the
#     new address as the link text ------------------------------------------
$ws3.Hyperlinks.Delete() | Out-Null
$ws3.Hyperlinks.Add(
    $ws3.Range("E2"),
    "mailto:interne-qualif@kobalt.fr",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "interne-qualif@kobalt.fr"
) | Out-Null

# Re-adding a hyperlink resets the cell to Excel's builtin blue/underlined
# "Hyperlink" style; put back the original look (plain blue Arial 10,
# no underline) that the rest of the workbook's e-mail cells use.
$ws3.Range("E2").Font.Name = "Arial"
$ws3.Range("E2").Font.Size = 10
$ws3.Range("E2").Font.Color = 16711680
$ws3.Range("E2").Font.Underline = $false

# --- Minor layout tweaks that came with the edit --------------------------
$ws3.Columns.Item(5).ColumnWidth = 20.18
$ws3.Rows.Item(2).RowHeight = 12.8

# --- TechnicalUser becomes the active/selected sheet (was BasicUser) ------
$ws3.Activate() | Out-Null
$ws3.Range("H2").Select() | Out-Null
